# ---------------------------------------------------------------------------
# "manual curve and PARAMS #6"
#
# Fills in the previously-blank replicate measurements (and their triplicate
# averages) for rows 6-17 of the "Growth_221117" sheet. Rows 2-5 already had
# data; rows 6-17 only had the time column (A) populated, with B/C/D and the
# F/G/H, J/K/L, N/O/P replicate triplets blank (or, for the latter, entirely
# absent) and the E/I/M/Q average formulas therefore evaluating to 0.
#
# Column E/I/M/Q already contain shared "=SUM(..)/3" formulas, so they (and
# the embedded chart that is fed straight from them) recompute automatically
# once the underlying replicate values below are entered - no formula edits
# are required there.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Growth_221117")
$ws.Activate()

# --- Row 6 (time = 2) ---
$ws.Range("B6").Value = 0.32
$ws.Range("C6").Value = 0.3
$ws.Range("D6").Value = 0.37
$ws.Range("F2:H2").Copy()
$ws.Range("F6:H6").PasteSpecial(-4122)
$ws.Range("F6").Value = 0.18
$ws.Range("G6").Value = 0.19
$ws.Range("H6").Value = 0.16
$ws.Range("J2:L2").Copy()
$ws.Range("J6:L6").PasteSpecial(-4122)
$ws.Range("J6").Value = 0.55
$ws.Range("K6").Value = 0.51
$ws.Range("L6").Value = 0.54
$ws.Range("N2:P2").Copy()
$ws.Range("N6:P6").PasteSpecial(-4122)
$ws.Range("N6").Value = 0.21
$ws.Range("O6").Value = 0.28
$ws.Range("P6").Value = 0.24

# --- Row 7 (time = 2.5) ---
$ws.Range("B7").Value = 0.47
$ws.Range("C7").Value = 0.46
$ws.Range("D7").Value = 0.49
$ws.Range("F2:H2").Copy()
$ws.Range("F7:H7").PasteSpecial(-4122)
$ws.Range("F7").Value = 0.25
$ws.Range("G7").Value = 0.24
$ws.Range("H7").Value = 0.25
$ws.Range("J2:L2").Copy()
$ws.Range("J7:L7").PasteSpecial(-4122)
$ws.Range("J7").Value = 0.89
$ws.Range("K7").Value = 0.86
$ws.Range("L7").Value = 0.89
$ws.Range("N2:P2").Copy()
$ws.Range("N7:P7").PasteSpecial(-4122)
$ws.Range("N7").Value = 0.3
$ws.Range("O7").Value = 0.25
$ws.Range("P7").Value = 0.33

# --- Row 8 (time = 3) ---
$ws.Range("B8").Value = 0.74
$ws.Range("C8").Value = 0.74
$ws.Range("D8").Value = 0.75
$ws.Range("F2:H2").Copy()
$ws.Range("F8:H8").PasteSpecial(-4122)
$ws.Range("F8").Value = 0.4
$ws.Range("G8").Value = 0.35
$ws.Range("H8").Value = 0.35
$ws.Range("J2:L2").Copy()
$ws.Range("J8:L8").PasteSpecial(-4122)
$ws.Range("J8").Value = 1.57
$ws.Range("K8").Value = 1.46
$ws.Range("L8").Value = 1.56
$ws.Range("N2:P2").Copy()
$ws.Range("N8:P8").PasteSpecial(-4122)
$ws.Range("N8").Value = 0.47
$ws.Range("O8").Value = 0.45
$ws.Range("P8").Value = 0.62

# --- Row 9 (time = 3.5) ---
$ws.Range("B9").Value = 1.17
$ws.Range("C9").Value = 1.02
$ws.Range("D9").Value = 1.01
$ws.Range("F2:H2").Copy()
$ws.Range("F9:H9").PasteSpecial(-4122)
$ws.Range("F9").Value = 0.51
$ws.Range("G9").Value = 0.46
$ws.Range("H9").Value = 0.54
$ws.Range("J2:L2").Copy()
$ws.Range("J9:L9").PasteSpecial(-4122)
$ws.Range("J9").Value = 2.71
$ws.Range("K9").Value = 2.35
$ws.Range("L9").Value = 2.4
$ws.Range("N2:P2").Copy()
$ws.Range("N9:P9").PasteSpecial(-4122)
$ws.Range("N9").Value = 0.77
$ws.Range("O9").Value = 0.76
$ws.Range("P9").Value = 0.87

# --- Row 10 (time = 4) ---
$ws.Range("B10").Value = 1.59
$ws.Range("C10").Value = 1.55
$ws.Range("D10").Value = 1.49
$ws.Range("F2:H2").Copy()
$ws.Range("F10:H10").PasteSpecial(-4122)
$ws.Range("F10").Value = 0.81
$ws.Range("G10").Value = 0.8
$ws.Range("H10").Value = 0.81
$ws.Range("J2:L2").Copy()
$ws.Range("J10:L10").PasteSpecial(-4122)
$ws.Range("J10").Value = 3.82
$ws.Range("K10").Value = 3.32
$ws.Range("L10").Value = 3.72
$ws.Range("N2:P2").Copy()
$ws.Range("N10:P10").PasteSpecial(-4122)
$ws.Range("N10").Value = 1.17
$ws.Range("O10").Value = 1.18
$ws.Range("P10").Value = 1.32

# --- Row 11 (time = 4.5) ---
$ws.Range("B11").Value = 2.02
$ws.Range("C11").Value = 2.04
$ws.Range("D11").Value = 2.01
$ws.Range("F2:H2").Copy()
$ws.Range("F11:H11").PasteSpecial(-4122)
$ws.Range("F11").Value = 1.17
$ws.Range("G11").Value = 1.16
$ws.Range("H11").Value = 1.99
$ws.Range("J2:L2").Copy()
$ws.Range("J11:L11").PasteSpecial(-4122)
$ws.Range("J11").Value = 5.36
$ws.Range("K11").Value = 4.69
$ws.Range("L11").Value = 4.91
$ws.Range("N2:P2").Copy()
$ws.Range("N11:P11").PasteSpecial(-4122)
$ws.Range("N11").Value = 1.65
$ws.Range("O11").Value = 1.64
$ws.Range("P11").Value = 1.71

# --- Row 12 (time = 5) ---
$ws.Range("B12").Value = 2.61
$ws.Range("C12").Value = 2.7
$ws.Range("D12").Value = 2.51
$ws.Range("F2:H2").Copy()
$ws.Range("F12:H12").PasteSpecial(-4122)
$ws.Range("F12").Value = 1.57
$ws.Range("G12").Value = 1.49
$ws.Range("H12").Value = 1.58
$ws.Range("J2:L2").Copy()
$ws.Range("J12:L12").PasteSpecial(-4122)
$ws.Range("J12").Value = 7.3
$ws.Range("K12").Value = 6.27
$ws.Range("L12").Value = 6.72
$ws.Range("N2:P2").Copy()
$ws.Range("N12:P12").PasteSpecial(-4122)
$ws.Range("N12").Value = 2.01
$ws.Range("O12").Value = 2.04
$ws.Range("P12").Value = 2.15

# --- Row 13 (time = 5.5) ---
$ws.Range("B13").Value = 3.26
$ws.Range("C13").Value = 3.25
$ws.Range("D13").Value = 3.21
$ws.Range("F2:H2").Copy()
$ws.Range("F13:H13").PasteSpecial(-4122)
$ws.Range("F13").Value = 1.98
$ws.Range("G13").Value = 1.85
$ws.Range("H13").Value = 1.94
$ws.Range("J2:L2").Copy()
$ws.Range("J13:L13").PasteSpecial(-4122)
$ws.Range("J13").Value = 8.74
$ws.Range("K13").Value = 7.86
$ws.Range("L13").Value = 8.62
$ws.Range("N2:P2").Copy()
$ws.Range("N13:P13").PasteSpecial(-4122)
$ws.Range("N13").Value = 2.62
$ws.Range("O13").Value = 2.56
$ws.Range("P13").Value = 2.68

# --- Row 14 (time = 6) ---
$ws.Range("B14").Value = 3.78
$ws.Range("C14").Value = 3.65
$ws.Range("D14").Value = 3.85
$ws.Range("F2:H2").Copy()
$ws.Range("F14:H14").PasteSpecial(-4122)
$ws.Range("F14").Value = 2.27
$ws.Range("G14").Value = 2.31
$ws.Range("H14").Value = 2.34
$ws.Range("J2:L2").Copy()
$ws.Range("J14:L14").PasteSpecial(-4122)
$ws.Range("J14").Value = 9.96
$ws.Range("K14").Value = 8.83
$ws.Range("L14").Value = 9.68
$ws.Range("N2:P2").Copy()
$ws.Range("N14:P14").PasteSpecial(-4122)
$ws.Range("N14").Value = 3.06
$ws.Range("O14").Value = 3.02
$ws.Range("P14").Value = 3.02

# --- Row 15 (time = 6.5) ---
# Recorded by the author as halves of a combined duplicate reading, so these
# (unlike every other data row) are literal "=x/2" formulas; and because they
# were typed straight into previously empty/absent cells the replicate-triplet
# columns (F/G/H, J/K/L, N/O/P) pick up the plain default cell style (no
# explicit style id), same as every other brand-new cell below.
$ws.Range("B15").Formula = "=9.36/2"
$ws.Range("C15").Value = 4.5
$ws.Range("D15").Formula = "=8.89/2"
$ws.Range("F15").Formula = "=4.97/2"
$ws.Range("G15").Formula = "=4.89/2"
$ws.Range("H15").Formula = "=4.91/2"
$ws.Range("J15").Formula = "=19.85/2"
$ws.Range("K15").Formula = "=18.37/2"
$ws.Range("L15").Formula = "=19.13/2"
$ws.Range("N15").Formula = "=6.61/2"
$ws.Range("O15").Formula = "=6.65/2"
$ws.Range("P15").Formula = "=6.62/2"

# --- Row 16 (time = 7) ---
$ws.Range("B16").Value = 5.57
$ws.Range("C16").Value = 5.86
$ws.Range("D16").Value = 5.64
$ws.Range("F2:H2").Copy()
$ws.Range("F16:H16").PasteSpecial(-4122)
$ws.Range("F16").Value = 3.19
$ws.Range("G16").Value = 3.2
$ws.Range("H16").Value = 3.23
$ws.Range("J2:L2").Copy()
$ws.Range("J16:L16").PasteSpecial(-4122)
$ws.Range("J16").Value = 11.99
$ws.Range("K16").Value = 10.84
$ws.Range("L16").Value = 11.47
$ws.Range("N2:P2").Copy()
$ws.Range("N16:P16").PasteSpecial(-4122)
$ws.Range("N16").Value = 3.92
$ws.Range("O16").Value = 3.96
$ws.Range("P16").Value = 4.04

# --- Row 17 (time = 7.5) ---
$ws.Range("B17").Value = 6.36
$ws.Range("C17").Value = 6.37
$ws.Range("D17").Value = 6.11
$ws.Range("F2:H2").Copy()
$ws.Range("F17:H17").PasteSpecial(-4122)
$ws.Range("F17").Value = 3.55
$ws.Range("G17").Value = 3.48
$ws.Range("H17").Value = 3.49
$ws.Range("J2:L2").Copy()
$ws.Range("J17:L17").PasteSpecial(-4122)
$ws.Range("J17").Value = 12.19
$ws.Range("K17").Value = 10.67
$ws.Range("L17").Value = 11.49
$ws.Range("N2:P2").Copy()
$ws.Range("N17:P17").PasteSpecial(-4122)
$ws.Range("N17").Value = 4.43
$ws.Range("O17").Value = 4.32
$ws.Range("P17").Value = 4.56

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("B18").Select()
